$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 681, shifting the existing rows 681:728 down to 682:729
$ws.Rows.Item(681).EntireRow.Insert()

# Populate the newly inserted row 681 with the new market record
$ws.Range("A681").Value = 5
$ws.Range("B681").Value = "Macroferia Regional de Talca"
$ws.Range("C681").Value = "Maule"
$ws.Range("D681").Value = 44931
$ws.Range("E681").Value = 7
$ws.Range("F681").Value = 100112002
$ws.Range("G681").Value = "Pimiento"
$ws.Range("H681").Value = "Cuatro cascos"
$ws.Range("I681").Value = "Primera"
$ws.Range("J681").Value = 400
$ws.Range("K681").Value = 8000
$ws.Range("L681").Value = 8000
$ws.Range("M681").Value = 8000
$ws.Range("N681").Value = "`$/caja 18 kilos"
$ws.Range("O681").Value = "Región del Maule"
$ws.Range("P681").Value = 444
$ws.Range("Q681").Value = 18
$ws.Range("R681").Value = "Hortaliza"
